$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert space for the 4 new control rows -------------------------------
# 1) Make room for "3.10" below the (still at row 10) "3.3" row, outside the
#    validated C4:C10 range so the data-validation sqref is left untouched.
$ws.Rows.Item(11).Insert()

# 2) Make room for "2.18" / "2.19" above the "3.3" row (still at row 10).
$ws.Range("A10:A11").EntireRow.Insert()

# 3) Make room for "1.22" above the "2.6" row (still at row 6).
$ws.Rows.Item(6).Insert()

# Control-ID cells ("1.22", "2.18", ...) look like numbers to Excel's type
# sniffer, so they'd otherwise be stored as floats (2.1800000000000002). Mark
# the cell as text first, write the value, then drop the now-unneeded "@"
# number format so the cell is left with the default (unstyled) look.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- New row 6: Control 1.22 ------------------------------------------------
Set-TextValue $ws.Range("A6") "1.22"
$ws.Range("B6").Value = "Information Barriers for AI Agents"
$ws.Range("C6").Value = "Not Started"
$ws.Range("D6:E6").Clear()

# --- New rows 11-12: Controls 2.18 / 2.19 -----------------------------------
Set-TextValue $ws.Range("A11") "2.18"
$ws.Range("B11").Value = "Automated Conflict of Interest Testing"
$ws.Range("C11").Value = "Not Started"
$ws.Range("D11:E11").Clear()

Set-TextValue $ws.Range("A12") "2.19"
$ws.Range("B12").Value = "Customer AI Disclosure and Transparency"
$ws.Range("C12").Value = "Not Started"
$ws.Range("D12:E12").Clear()

# --- New row 14: Control 3.10 -----------------------------------------------
Set-TextValue $ws.Range("A14") "3.10"
$ws.Range("B14").Value = "Hallucination Feedback Loop"
$ws.Range("C14").Value = "Not Started"
$ws.Range("D14:E14").Clear()

# --- Restore the data-validation sqref (row inserts inside C4:C10 grow it) --
$v = $ws.Range("C4:C13").Validation
$v.Delete()
$v2 = $ws.Range("C4:C10").Validation
$v2.Add(3, 1, 1, """Not Started,In Progress,Completed,N/A""")
$v2.ShowInput = $False
$v2.ShowError = $False

# --- Materialize the blank spacer rows (row 2, and rows 15-16 before the
#     footer) as empty row records, matching the source sheet's style of
#     using placeholder rows for vertical spacing. Setting OutlineLevel to
#     its already-default value (0) touches the row without adding any
#     cell/format content.
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Rows.Item(15).OutlineLevel = 0
$ws.Rows.Item(16).OutlineLevel = 0
